$d = $word.ActiveDocument

# The site generator rebuild dropped the trailing "Ver no Jupiter..." line and
# the "(c) 2020 ..." Jekyll footer line (plus the blank spacer paragraph that
# preceded them) that used to be appended right after the "Requisitos"
# section's last entry (the "LOQ4073..." paragraph).
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*LOQ4073*") {
        # The next three paragraphs are: the blank spacer, the
        # "Ver no Jupiter Salvar em pdf Salvar em docx" line, and the
        # "© 2020 . Contact: ..." footer line. Remove all three, including
        # their paragraph marks, in one shot.
        $startPara = $d.Paragraphs.Item($i + 1)
        $endPara = $d.Paragraphs.Item($i + 3)
        $r = $d.Range($startPara.Range.Start, $endPara.Range.End)
        $r.Delete()
        break
    }
}
